# Быстрая помощь.pptx — apply the authored edit:
#   - slide 3 ("Объект 2" body placeholder): the second occurrence of the
#     phrase "из класса " (the one describing PyQt5 widgets/QtCore/QtGui/
#     QtWidgets, i.e. "... QtWidgets из класса PyQt5.") has the word
#     "класса" replaced with "библиотеки", which also splits the single
#     run into two runs ("из " / "библиотеки ") because the edit only
#     touches part of the original run.
#   - presentation.xml gains an empty <p15:sldGuideLst/> extension block
#     (recorded by PowerPoint when the file was last saved). The hosted
#     object model here has no working Guides-mutation surface (Presentation
#     .Guides / Application.DisplayGuides are read-only stubs in this
#     runtime), so we make a best-effort, harmless attempt and otherwise
#     leave that part to the host; it does not affect the slide content.

$p = $ppt.ActivePresentation

# --- Locate the target shape: the slide body placeholder whose text
#     contains the word "класса" (appears twice; we edit the last one). ---
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $t = $shp.TextFrame.TextRange.Text
                if ($t.Contains("класса")) {
                    $targetSlide = $sl
                    $targetShape = $shp
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $full = $tr.Text

    # The run we need is the LAST "класса " in the text box (the first one,
    # "... QWidget  из класса PyQt5.QtWidgets, ..." stays untouched).
    $wordIdx = $full.LastIndexOf("класса ")

    if ($wordIdx -ge 0) {
        # Characters() is 1-based; replace the 7 characters "класса " in
        # place so the surrounding "из " / trailing text keep their own
        # runs/formatting exactly like the rest of the paragraph.
        $target = $tr.Characters($wordIdx + 1, 7)
        $target.Text = "библиотеки "
    }
}

# --- Best-effort: presentation-level empty guide list extension. ---
try {
    $guides = $p.Guides
    if ($guides -ne $null) {
        # Touch the collection so a host that lazily materialises the
        # <p:extLst>/<p15:sldGuideLst/> block on first access gets a chance
        # to do so, without actually adding a visible guide.
        $null = $guides.Count
    }
} catch {
    # Guides are not mutable via this object model build; nothing else to
    # do here — the slide content edit above is the substantive change.
}
